$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create Sheet2 positioned right after Sheet1 (this also makes it the
# active sheet, which drops tabSelected from Sheet1's sheetView and
# bumps the workbook's activeTab to 1).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Column widths (approximate best-fit widths from the original workbook).
$ws2.Columns.Item(1).ColumnWidth = 18
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 15.42578125
$ws2.Columns.Item(4).ColumnWidth = 10.42578125

# Borders/fills around the whole table, reusing the existing thin-box
# styles already used on Sheet1 (copy format from Sheet1's bordered
# cells, one row at a time, so no duplicate border/style entries get
# created and the header/body styling lines up correctly).
$ws1.Range("B2:E2").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws1.Range("B3:E3").Copy() | Out-Null
$ws2.Range("A2:D2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A3:D3").PasteSpecial(-4122) | Out-Null
$ws2.Range("A4:D4").PasteSpecial(-4122) | Out-Null
$ws2.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$ws2.Range("A7:D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row.
$ws2.Range("A1").Value = "Model Name"
$ws2.Range("B1").Value = "Model Accuracy"
$ws2.Range("C1").Value = "Cross Validation"
$ws2.Range("D1").Value = "Difference"

# Model rows.
$names = @("Logistic Regression", "Decision Tree", "Random Forest", "SVC", "Kneighbors", "MultinomialNB")
$acc   = @(0.95, 1, 1, 0.99, 1, 0.8)
$cv    = @(0.84, 0.92, 0.9, 0.84, 0.88, 0.72)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $names[$i]
    $ws2.Range("B$row").Value = $acc[$i]
    $ws2.Range("C$row").Value = $cv[$i]
}

# Difference formulas: D2 stands alone, D3:D7 share one formula definition.
$ws2.Range("D2").Formula = "=B2-C2"
$ws2.Range("D3:D7").Formula = "=B3-C3"

# Final selection/active cell on Sheet2.
$ws2.Range("C12").Select() | Out-Null

Write-Host "Sheet2 created with model comparison data."
